# Case and Fatality Demographics Data Updated
# Updates the "Fatalities by Age Group", "Fatalities by Gender" and
# "Fatalities by Race-Ethnicity" sheets with the 9/02/21 report numbers,
# and moves the active/selected sheet from "Fatalities by Age Group" to
# "Fatalities by Race-Ethnicity".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Fatalities by Age Group
# ---------------------------------------------------------------------
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")

$wsAge.Range("B2").Value = 10
$wsAge.Range("B4").Value = 43
$wsAge.Range("B5").Value = 368
$wsAge.Range("B6").Value = 1177
$wsAge.Range("B7").Value = 3226
$wsAge.Range("B8").Value = 6801
$wsAge.Range("B9").Value = 5420
$wsAge.Range("B10").Value = 6804
$wsAge.Range("B11").Value = 7472
$wsAge.Range("B12").Value = 7282
$wsAge.Range("B13").Value = 17981

# ---------------------------------------------------------------------
# Fatalities by Gender
# ---------------------------------------------------------------------
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")

$wsGender.Range("B2").Value = 23678
$wsGender.Range("B3").Value = 32919

# ---------------------------------------------------------------------
# Fatalities by Race-Ethnicity
# ---------------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")

$wsRace.Range("B2").Value = 1144
$wsRace.Range("B3").Value = 5753
$wsRace.Range("B4").Value = 25727
$wsRace.Range("B5").Value = 316
$wsRace.Range("B6").Value = 23631
$wsRace.Range("B7").Value = 27

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping, matching the saved workbook view
# ---------------------------------------------------------------------
[void]$wsAge.Range("G8").Select()

[void]$wsRace.Activate()
[void]$wsRace.Range("B21").Select()
